# Commit: "Have shown to Shachi ma'am regarding excel of Automation"
#
# - Renames the "DashboardData" sheet to "CategoryData" and turns it into a
#   Category-import template (header row: Category Name / Category Image /
#   Test Results).
# - Clears the stale "Test Results" column (C) on the LoginData sheet - those
#   values were leftover Selenium run noise (long alert/stack-trace strings).
# - Touches up column widths / active selection on both sheets to match what
#   was left on screen.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Rename sheet 2: DashboardData -> CategoryData
# ---------------------------------------------------------------------------
$wsLogin = $wb.Worksheets.Item("LoginData")
$wsCat   = $wb.Worksheets.Item("DashboardData")
$wsCat.Name = "CategoryData"

# ---------------------------------------------------------------------------
# 2. LoginData: clear the "Test Results" column (C2:C5) - the automation-run
#    noise (pass/fail alert dumps) is no longer needed.
# ---------------------------------------------------------------------------
$wsLogin.Range("C2:C5").ClearContents()

# Column widths widened/narrowed slightly after the cleanup.
$wsLogin.Columns.Item(1).ColumnWidth = 21.93
$wsLogin.Columns.Item(2).ColumnWidth = 24.93
$wsLogin.Columns.Item(3).ColumnWidth = 46.22
$wsLogin.Columns.Item(4).ColumnWidth = 7.08
$wsLogin.Columns.Item(5).Resize(1, 1021).EntireColumn.ColumnWidth = -0.5

# Active selection left on C2:C5 after clearing it.
$wsLogin.Activate() | Out-Null
$wsLogin.Range("C2:C5").Select() | Out-Null

# ---------------------------------------------------------------------------
# 3. CategoryData: add the header row and size the sheet like the new data.
# ---------------------------------------------------------------------------
$wsLogin.Range("A1:C1").Copy()
$wsCat.Range("A1:C1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$wsCat.Range("A1").Value = "Category Name "
$wsCat.Range("B1").Value = "Category Image"
$wsCat.Range("C1").Value = "Test Results"
$wsCat.Rows.Item(1).RowHeight = 30

$wsCat.Columns.Item(1).ColumnWidth = 26.51
$wsCat.Columns.Item(2).ColumnWidth = 35.79
$wsCat.Columns.Item(3).ColumnWidth = 24.36
$wsCat.Columns.Item(4).Resize(1, 1022).EntireColumn.ColumnWidth = -0.21

$wsCat.Activate() | Out-Null
$u = $excel.Union($wsCat.Range("C2:C5"), $wsCat.Range("A3"))
$u.Select() | Out-Null
